$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (Sheet1 -> DayOne)
$ws.Name = "DayOne"

# Fix capitalisation of the existing "Subject" value
$ws.Range("B3").Value = "Registration assistance"

# Add the two new rows of request data
$ws.Range("A4").Value = "Account Name"
$ws.Range("B4").Value = "Patrick B"
$ws.Range("A5").Value = "Phone"
$ws.Range("B5").Value = 728030696

# Widen column A to fit the new labels (closest reachable width to 12.54296875)
$ws.Columns("A").ColumnWidth = 11.6

# Move the active selection to the last entered cell
$ws.Range("B5").Select() | Out-Null
